$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-absent rows 18-23 (table shrinks from 22 to 16 data rows)
$ws.Rows("18:23").Delete()

# Helper: write a value as text (even if it looks numeric), with no residual style
function Set-TextCell($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$ws.Cells.Item(2,1).Value = "Queneau_ref vs Queneau_gen"
$ws.Cells.Item(2,2).Value = "Structural"
$ws.Cells.Item(2,3).Value = 4.193026965503968
Set-TextCell $ws.Cells.Item(2,4) "0.000071712258"
$ws.Cells.Item(2,5).Value = "Highly Significant"
$ws.Cells.Item(2,6).Value = 1.869331015836096

$ws.Cells.Item(3,1).Value = "Queneau_ref vs Queneau_gen"
$ws.Cells.Item(3,2).Value = "Indexes"
$ws.Cells.Item(3,3).Value = -18.15582301780313
Set-TextCell $ws.Cells.Item(3,4) "0.000000000000"
$ws.Cells.Item(3,5).Value = "Highly Significant"
$ws.Cells.Item(3,6).Value = -638.5826499449554

$ws.Cells.Item(4,1).Value = "Queneau_ref vs Queneau_gen"
$ws.Cells.Item(4,2).Value = "Letters"
$ws.Cells.Item(4,3).Value = -5.32126071122836
Set-TextCell $ws.Cells.Item(4,4) "0.000000451835"
$ws.Cells.Item(4,5).Value = "Highly Significant"
$ws.Cells.Item(4,6).Value = -0.001454911388067987

$ws.Cells.Item(5,1).Value = "Queneau_ref vs Queneau_gen"
$ws.Cells.Item(5,2).Value = "Punctuation"
$ws.Cells.Item(5,3).Value = -5.327845847264404
Set-TextCell $ws.Cells.Item(5,4) "0.000000374761"
$ws.Cells.Item(5,5).Value = "Highly Significant"
$ws.Cells.Item(5,6).Value = -0.001038611566991879

# Row 6 is unchanged in the target, left untouched

$ws.Cells.Item(7,1).Value = "Queneau_ref vs Queneau_gen"
$ws.Cells.Item(7,2).Value = "NER"
$ws.Cells.Item(7,3).Value = 2.472802896851884
Set-TextCell $ws.Cells.Item(7,4) "0.015195208041"
$ws.Cells.Item(7,5).Value = "Significant"
$ws.Cells.Item(7,6).Value = 0.177497317873082

$ws.Cells.Item(8,1).Value = "Queneau_ref vs Queneau_gen"
$ws.Cells.Item(8,2).Value = "Function words"
$ws.Cells.Item(8,3).Value = -5.118552802480382
Set-TextCell $ws.Cells.Item(8,4) "0.000001240861"
$ws.Cells.Item(8,5).Value = "Highly Significant"
$ws.Cells.Item(8,6).Value = -0.0004868362043281909

$ws.Cells.Item(9,1).Value = "Queneau_ref vs Queneau_gen"
$ws.Cells.Item(9,2).Value = "Numbers"
$ws.Cells.Item(9,3).Value = 1.328138461039664
Set-TextCell $ws.Cells.Item(9,4) "0.186400714163"
$ws.Cells.Item(9,5).Value = "Not Significant"
$ws.Cells.Item(9,6).Value = 0.006183489104453266

$ws.Cells.Item(10,1).Value = "Feneon_ref vs Queneau_gen"
$ws.Cells.Item(10,2).Value = "Structural"
$ws.Cells.Item(10,3).Value = 4.147606195723885
Set-TextCell $ws.Cells.Item(10,4) "0.000057193884"
$ws.Cells.Item(10,5).Value = "Highly Significant"
$ws.Cells.Item(10,6).Value = 0.5339516746388488

$ws.Cells.Item(11,1).Value = "Feneon_ref vs Queneau_gen"
$ws.Cells.Item(11,2).Value = "Indexes"
$ws.Cells.Item(11,3).Value = 1.69954272643086
Set-TextCell $ws.Cells.Item(11,4) "0.091455763127"
$ws.Cells.Item(11,5).Value = "Not Significant"
$ws.Cells.Item(11,6).Value = 67.61067444662422

$ws.Cells.Item(12,1).Value = "Feneon_ref vs Queneau_gen"
$ws.Cells.Item(12,2).Value = "Letters"
$ws.Cells.Item(12,3).Value = -6.078157073799479
Set-TextCell $ws.Cells.Item(12,4) "0.000000010433"
$ws.Cells.Item(12,5).Value = "Highly Significant"
$ws.Cells.Item(12,6).Value = -0.001883426952475654

$ws.Cells.Item(13,1).Value = "Feneon_ref vs Queneau_gen"
$ws.Cells.Item(13,2).Value = "Punctuation"
$ws.Cells.Item(13,3).Value = -0.6881206213851467
Set-TextCell $ws.Cells.Item(13,4) "0.492562470693"
$ws.Cells.Item(13,5).Value = "Not Significant"
$ws.Cells.Item(13,6).Value = -0.000155805011576084

$ws.Cells.Item(14,1).Value = "Feneon_ref vs Queneau_gen"
$ws.Cells.Item(14,2).Value = "TAG"
$ws.Cells.Item(14,3).Value = 7.053747288129914
Set-TextCell $ws.Cells.Item(14,4) "0.000000000075"
$ws.Cells.Item(14,5).Value = "Highly Significant"
$ws.Cells.Item(14,6).Value = 0.267419245729748

$ws.Cells.Item(15,1).Value = "Feneon_ref vs Queneau_gen"
$ws.Cells.Item(15,2).Value = "NER"
$ws.Cells.Item(15,3).Value = 5.198731596449358
Set-TextCell $ws.Cells.Item(15,4) "0.000000818268"
$ws.Cells.Item(15,5).Value = "Highly Significant"
$ws.Cells.Item(15,6).Value = 0.2579908675799087

$ws.Cells.Item(16,1).Value = "Feneon_ref vs Queneau_gen"
$ws.Cells.Item(16,2).Value = "Function words"
$ws.Cells.Item(16,3).Value = -8.086454468937944
Set-TextCell $ws.Cells.Item(16,4) "0.000000000000"
$ws.Cells.Item(16,5).Value = "Highly Significant"
$ws.Cells.Item(16,6).Value = -0.0009343893458050375

$ws.Cells.Item(17,1).Value = "Feneon_ref vs Queneau_gen"
$ws.Cells.Item(17,2).Value = "Numbers"
$ws.Cells.Item(17,3).Value = 2.011576279521369
Set-TextCell $ws.Cells.Item(17,4) "0.046408495312"
$ws.Cells.Item(17,5).Value = "Significant"
$ws.Cells.Item(17,6).Value = 0.01012054867125877

